# Applies the edit described by the commit:
# "se descargan los logs cuando se crga el udas salga o no error"
#
# Functional changes to the workbook's lookup/list sheets:
#  - hardware: remove the duplicated "CELULAR SAMSUNG SM-J510MN" row
#  - microphone: populate the (previously empty) list of microphone types
#  - gain: populate the (previously empty) list of gain settings
#  - funding: append a new funding source, "ECOPETROL"

$wb = $excel.ActiveWorkbook

# --- hardware: drop the duplicate "CELULAR SAMSUNG SM-J510MN" row ---------
$ws = $wb.Worksheets.Item("hardware")
$ws.Rows.Item(12).Delete()

# --- microphone: fill in the microphone-type list --------------------------
$ws = $wb.Worksheets.Item("microphone")
$values = @(
    "NO SE CONOCE",
    "MICROFONO INTEGRADO",
    "AUDIOMOTH EXTERNO",
    "SMX-II",
    "SENNHEISER DIRECCIONAL"
)
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# --- gain: fill in the gain-setting list -----------------------------------
$ws = $wb.Worksheets.Item("gain")
$values = @(
    "NO SE CONOCE",
    "AUDIOMOTH GANANCIA BAJA",
    "AUDIOMOTH GANANCIA BAJA-MEDIA",
    "AUDIOMOTH GANANCIA MEDIA",
    "AUDIOMOTH GANANCIA MEDIA-ALTA",
    "AUDIOMOTH GANANCIA ALTA",
    "16DB",
    "SM4 SIN GANANCIA",
    "AUDIOMOTH SIN GANANCIA",
    "AUDIOMOTH SIN GANANCIA 2",
    "AUDIOMOTH SIN GANANCIA 3"
)
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# --- funding: append "ECOPETROL" as a new funding source -------------------
$ws = $wb.Worksheets.Item("funding")
$ws.Cells.Item(8, 1).Value = "ECOPETROL"
